$wb = $excel.ActiveWorkbook

# --- mobility sheet ---
$ws = $wb.Worksheets.Item("mobility")
$ws.Activate()

# Update values: regroup_time 0.01 -> 3, min_dist_to_derive 750000 -> 700000
$ws.Range("B20").Value = 3
$ws.Range("B24").Value = 700000

# Update the selected/active cell to B8
$ws.Range("B8").Select()

# --- infrastructure sheet ---
$ws2 = $wb.Worksheets.Item("infrastructure")
$ws2.Activate()
$ws2.Range("B1").Select()

# --- categories sheet ---
$ws4 = $wb.Worksheets.Item("categories")
$ws4.Activate()
$ws4.Range("B8").Select()

# Restore original active sheet (mobility was tabSelected in the source file)
$ws.Activate()
